# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
# Change cell B11 on the active sheet from "R40" to the text value "1".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = "1"
